$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace column B (previously "Value") with the previously-computed
# "% Change vs Last Year" values that used to live in column C, before
# column C (and its formulas) is removed.
$ws.Range("B2").Value = 3.1522570772761993
$ws.Range("B3").Value = 2.8820220389920381
$ws.Range("B4").Value = 3.4719513146602266
$ws.Range("B5").Value = 2.6597825236642336
$ws.Range("B6").Value = 3.2303925440328563
$ws.Range("B7").Value = 3.0247697681803798
$ws.Range("B8").Value = 3.5131642707155519
$ws.Range("B9").Value = 6.8366067697450905
$ws.Range("B10").Value = 6.5740740740740655
$ws.Range("B11").Value = 7.5110959371799124
$ws.Range("B12").Value = 7.5947167188042952
$ws.Range("B13").Value = 2.9158782040254572
$ws.Range("B14").Value = 1.2615069894306163
$ws.Range("B15").Value = -11.309613928841788
$ws.Range("B16").Value = 1.8226862502211949
$ws.Range("B17").Value = 8.128720238095255
$ws.Range("B18").Value = 10.689687706387385
$ws.Range("B19").Value = 25.977493801258824
$ws.Range("B20").Value = 7.8435114503816683
$ws.Range("B21").Value = 3.7237121358286718
$ws.Range("B22").Value = 3.2035053554040793
$ws.Range("B23").Value = 2.1429962984609308
$ws.Range("B24").Value = 2.0547278215989895
$ws.Range("B25").Value = 1.727183513248276
$ws.Range("B26").Value = 2.7102710271027286
$ws.Range("B27").Value = -6.4430875786020252
$ws.Range("B28").Value = -6.6963474468471729
$ws.Range("B29").Value = -6.3160798014158193
$ws.Range("B30").Value = -7.6986984214898913
$ws.Range("B31").Value = 1.6018518518518654
$ws.Range("B32").Value = 2.6487595597836311
$ws.Range("B33").Value = -0.4211297262656899
$ws.Range("B34").Value = -0.24861878453038555
$ws.Range("B35").Value = -0.96286107290233236
$ws.Range("B36").Value = -1.750206176120217
$ws.Range("B37").Value = 1.3359309768995242
$ws.Range("B38").Value = 0.78886310904870971
$ws.Range("B39").Value = 2.2887158803114094
$ws.Range("B40").Value = 4.2510508215513942
$ws.Range("B41").Value = 3.6043829296424468
$ws.Range("B42").Value = 4.5710403726707982
$ws.Range("B43").Value = 4.3661282427802206
$ws.Range("B44").Value = 2.8796068796068885
$ws.Range("B45").Value = 2.7961663867206887
$ws.Range("B46").Value = 2.7215631542219176
$ws.Range("B47").Value = -9.2161393530039053

# Rows 48-51 never had a computed % change (the shared formula only
# covered C2:C47), so their old column-B "Value" numbers are simply dropped.
$ws.Range("B48").ClearContents()
$ws.Range("B49").ClearContents()
$ws.Range("B50").ClearContents()
$ws.Range("B51").ClearContents()

# Drop the now-redundant column C (formulas + header) entirely; this
# shifts nothing since B already holds the right data, it just removes C.
$ws.Columns("C").Delete()

# Column B header becomes the old column-C header text.
$ws.Range("B1").Value = "% Change vs Last Year"
